$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3759.4546
$ws.Range("I62").Value = 3901.7058
$ws.Range("J62").Value = 3275.8
$ws.Range("K62").Value = 3901.7058
$ws.Range("L62").Value = 3275.8
$ws.Range("M62").Value = -3277.7058
$ws.Range("N62").Value = -4523.8

$ws.Range("H65").Value = 3759.4546
$ws.Range("I65").Value = 3901.7058
$ws.Range("J65").Value = 3275.8
$ws.Range("K65").Value = 19508.529
$ws.Range("L65").Value = 16379
$ws.Range("M65").Value = -16388.529
$ws.Range("N65").Value = -22619

$ws.Range("H101").Value = 31480.166
$ws.Range("I101").Value = 681.5
$ws.Range("J101").Value = 93077.5
$ws.Range("K101").Value = 2044.5
$ws.Range("L101").Value = 279232.5
$ws.Range("M101").Value = -422.5
$ws.Range("N101").Value = -282476.5

$ws.Range("H106").Value = 21619.85
$ws.Range("I106").Value = 24226.133
$ws.Range("J106").Value = 13801
$ws.Range("K106").Value = 24226.133
$ws.Range("L106").Value = 13801
$ws.Range("M106").Value = -23595.133
$ws.Range("N106").Value = -15063

$ws.Range("H113").Value = 15634796
$ws.Range("I113").Value = 41677500
$ws.Range("K113").Value = 41677500
$ws.Range("M113").Value = -41674246

$ws.Range("H137").Value = 1300.7241
$ws.Range("I137").Value = 1336.8572
$ws.Range("K137").Value = 4010.5716
$ws.Range("M137").Value = -1460.5716

$ws.Range("H138").Value = 4561.6055
$ws.Range("I138").Value = 1337.381
$ws.Range("J138").Value = 5792.673
$ws.Range("K138").Value = 4012.143
$ws.Range("L138").Value = 17378.019
$ws.Range("M138").Value = 1127.857
$ws.Range("N138").Value = -27658.019

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2890
$ws.Range("I32").Value = 2491.2769
$ws.Range("K32").Value = 2491.2769
$ws.Range("M32").Value = -2204.2769

$ws.Range("H61").Value = 2671.913
$ws.Range("I61").Value = 2516.147
$ws.Range("K61").Value = 2516.147
$ws.Range("M61").Value = -2304.147

$ws.Range("H74").Value = 2323.3547
$ws.Range("I74").Value = 1968.9131
$ws.Range("K74").Value = 1968.9131
$ws.Range("M74").Value = -1094.9131

$ws.Range("H77").Value = 2323.3547
$ws.Range("I77").Value = 1968.9131
$ws.Range("K77").Value = 9844.565500000001
$ws.Range("M77").Value = -5476.565500000001

$ws.Range("H136").Value = 2671.913
$ws.Range("I136").Value = 2516.147
$ws.Range("K136").Value = 7548.441
$ws.Range("M136").Value = -4998.441

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2222.6
$ws.Range("I105").Value = 1798.1428
$ws.Range("K105").Value = 1798.1428
$ws.Range("M105").Value = -51.14280000000008

$ws.Range("H107").Value = 5617.1035
$ws.Range("I107").Value = 4974.8096
$ws.Range("K107").Value = 4974.8096
$ws.Range("M107").Value = -3054.8096

$ws.Range("H134").Value = 41669836
$ws.Range("I134").Value = 45457820
$ws.Range("K134").Value = 136373460
$ws.Range("M134").Value = -136370925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1510.4286
$ws.Range("I31").Value = 1668.25
$ws.Range("K31").Value = 1668.25
$ws.Range("M31").Value = -1373.25

$ws.Range("H34").Value = 1510.4286
$ws.Range("I34").Value = 1668.25
$ws.Range("K34").Value = 1668.25
$ws.Range("M34").Value = -1466.25

$ws.Range("H39").Value = 17407
$ws.Range("I39").Value = 17407
$ws.Range("K39").Value = 17407
$ws.Range("M39").Value = -17016

$ws.Range("H49").Value = 17407
$ws.Range("I49").Value = 17407
$ws.Range("K49").Value = 17407
$ws.Range("M49").Value = -17225

$ws.Range("H99").Value = 3008.2727
$ws.Range("I99").Value = 2448.5
$ws.Range("J99").Value = 3328.1428
$ws.Range("K99").Value = 2448.5
$ws.Range("L99").Value = 3328.1428
$ws.Range("M99").Value = -950.5
$ws.Range("N99").Value = -6324.1428

$ws.Range("H122").Value = 3290.9285
$ws.Range("I122").Value = 2014
$ws.Range("K122").Value = 6042
$ws.Range("M122").Value = -3592

$ws.Range("H126").Value = 3008.2727
$ws.Range("I126").Value = 2448.5
$ws.Range("J126").Value = 3328.1428
$ws.Range("K126").Value = 7345.5
$ws.Range("L126").Value = 9984.428400000001
$ws.Range("M126").Value = -4875.5
$ws.Range("N126").Value = -14924.4284

$ws.Range("H134").Value = 1207.8064
$ws.Range("I134").Value = 1133.2222
$ws.Range("K134").Value = 3399.6666
$ws.Range("M134").Value = -864.6665999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1096.3334
$ws.Range("I47").Value = 994.5
$ws.Range("J47").Value = 1300
$ws.Range("K47").Value = 2983.5
$ws.Range("L47").Value = 3900
$ws.Range("M47").Value = -2552.5
$ws.Range("N47").Value = -4762

$ws.Range("H134").Value = 3110.6843
$ws.Range("I134").Value = 2733.5
$ws.Range("K134").Value = 8200.5
$ws.Range("M134").Value = -3130.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 999
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 999
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 999
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -1231

$ws.Range("H14").Value = 3424.6
$ws.Range("I14").Value = 5048
$ws.Range("K14").Value = 5048
$ws.Range("M14").Value = -4880

$ws.Range("H102").Value = 3427.2856
$ws.Range("I102").Value = 2998.25
$ws.Range("J102").Value = 3999.3333
$ws.Range("K102").Value = 2998.25
$ws.Range("L102").Value = 3999.3333
$ws.Range("M102").Value = -1376.25
$ws.Range("N102").Value = -7243.3333

$ws.Range("H107").Value = 52632704
$ws.Range("J107").Value = 742.2222
$ws.Range("L107").Value = 742.2222
$ws.Range("N107").Value = -4582.2222

$ws.Range("H113").Value = 2425.3076
$ws.Range("I113").Value = 2292.7144
$ws.Range("J113").Value = 2580
$ws.Range("K113").Value = 2292.7144
$ws.Range("L113").Value = 2580
$ws.Range("M113").Value = -122.7143999999998
$ws.Range("N113").Value = -6920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5713.524
$ws.Range("I40").Value = 5610.5
$ws.Range("J40").Value = 6043.2
$ws.Range("K40").Value = 5610.5
$ws.Range("L40").Value = 6043.2
$ws.Range("M40").Value = -5474.5
$ws.Range("N40").Value = -6315.2

$ws.Range("H122").Value = 4391.185
$ws.Range("I122").Value = 4154.7393
$ws.Range("K122").Value = 12464.2179
$ws.Range("M122").Value = -10014.2179

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 29666.666
$ws.Range("I32").Value = 29666.666
$ws.Range("K32").Value = 29666.666
$ws.Range("M32").Value = -29349.666
